$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Swap / update existing rows (B..AC only; A and E unchanged) ----
# Row 90
$ws.Cells.Item(90,2).Value = 7480684
$ws.Cells.Item(90,6).Value = "NC Magra"
$ws.Cells.Item(90,7).Value = "MC El Bayadh"
$ws.Cells.Item(90,8).Value = 1
$ws.Cells.Item(90,9).Value = 1
$ws.Cells.Item(90,10).Value = "D"
$ws.Cells.Item(90,11).Value = 1.85
$ws.Cells.Item(90,12).Value = 3
$ws.Cells.Item(90,13).Value = 4.25
$ws.Cells.Item(90,14).Value = 2
$ws.Cells.Item(90,15).Value = 2.9
$ws.Cells.Item(90,16).Value = 3.8
$ws.Cells.Item(90,17).Value = -0.25
$ws.Cells.Item(90,18).Value = 1.75
$ws.Cells.Item(90,19).Value = 2.05
$ws.Cells.Item(90,20).Value = 1.75
$ws.Cells.Item(90,21).Value = 1.75
$ws.Cells.Item(90,22).Value = 2.05
$ws.Cells.Item(90,23).Value = -1
$ws.Cells.Item(90,24).Value = 1.9
$ws.Cells.Item(90,25).Value = -1
$ws.Cells.Item(90,26).Value = -0.5
$ws.Cells.Item(90,27).Value = 0.5249999999999999
$ws.Cells.Item(90,28).Value = 0.375
$ws.Cells.Item(90,29).Value = -0.5

# Row 91
$ws.Cells.Item(91,2).Value = 7480688
$ws.Cells.Item(91,6).Value = "Paradou AC"
$ws.Cells.Item(91,7).Value = "CS Constantine"
$ws.Cells.Item(91,8).Value = 0
$ws.Cells.Item(91,9).Value = 0
$ws.Cells.Item(91,10).Value = "D"
$ws.Cells.Item(91,11).Value = 1.7
$ws.Cells.Item(91,12).Value = 3.2
$ws.Cells.Item(91,13).Value = 4.75
$ws.Cells.Item(91,14).Value = 1.95
$ws.Cells.Item(91,15).Value = 3
$ws.Cells.Item(91,16).Value = 3.8
$ws.Cells.Item(91,17).Value = -0.5
$ws.Cells.Item(91,18).Value = 1.975
$ws.Cells.Item(91,19).Value = 1.825
$ws.Cells.Item(91,20).Value = 2
$ws.Cells.Item(91,21).Value = 1.925
$ws.Cells.Item(91,22).Value = 1.875
$ws.Cells.Item(91,23).Value = -1
$ws.Cells.Item(91,24).Value = 2
$ws.Cells.Item(91,25).Value = -1
$ws.Cells.Item(91,26).Value = -1
$ws.Cells.Item(91,27).Value = 0.825
$ws.Cells.Item(91,28).Value = -1
$ws.Cells.Item(91,29).Value = 0.875

# Row 122
$ws.Cells.Item(122,2).Value = 7624657
$ws.Cells.Item(122,6).Value = "US Souf"
$ws.Cells.Item(122,7).Value = "Paradou AC"
$ws.Cells.Item(122,8).Value = 1
$ws.Cells.Item(122,9).Value = 4
$ws.Cells.Item(122,10).Value = "A"
$ws.Cells.Item(122,11).Value = 3
$ws.Cells.Item(122,12).Value = 2.8
$ws.Cells.Item(122,13).Value = 2.375
$ws.Cells.Item(122,14).Value = 4.333
$ws.Cells.Item(122,15).Value = 3
$ws.Cells.Item(122,16).Value = 1.85
$ws.Cells.Item(122,17).Value = 0.5
$ws.Cells.Item(122,18).Value = 1.85
$ws.Cells.Item(122,19).Value = 1.95
$ws.Cells.Item(122,20).Value = 2
$ws.Cells.Item(122,21).Value = 2.025
$ws.Cells.Item(122,22).Value = 1.775
$ws.Cells.Item(122,23).Value = -1
$ws.Cells.Item(122,24).Value = -1
$ws.Cells.Item(122,25).Value = 0.8500000000000001
$ws.Cells.Item(122,26).Value = -1
$ws.Cells.Item(122,27).Value = 0.95
$ws.Cells.Item(122,28).Value = 1.025
$ws.Cells.Item(122,29).Value = -1

# Row 124
$ws.Cells.Item(124,2).Value = 7656900
$ws.Cells.Item(124,6).Value = "NC Magra"
$ws.Cells.Item(124,7).Value = "ES Setif"
$ws.Cells.Item(124,8).Value = 0
$ws.Cells.Item(124,9).Value = 1
$ws.Cells.Item(124,10).Value = "A"
$ws.Cells.Item(124,11).Value = 1.833
$ws.Cells.Item(124,12).Value = 3.1
$ws.Cells.Item(124,13).Value = 4
$ws.Cells.Item(124,14).Value = 2.875
$ws.Cells.Item(124,15).Value = 3
$ws.Cells.Item(124,16).Value = 2.625
$ws.Cells.Item(124,17).Value = 0
$ws.Cells.Item(124,18).Value = 1.975
$ws.Cells.Item(124,19).Value = 1.825
$ws.Cells.Item(124,20).Value = 2
$ws.Cells.Item(124,21).Value = 1.925
$ws.Cells.Item(124,22).Value = 1.875
$ws.Cells.Item(124,23).Value = -1
$ws.Cells.Item(124,24).Value = -1
$ws.Cells.Item(124,25).Value = 1.625
$ws.Cells.Item(124,26).Value = -1
$ws.Cells.Item(124,27).Value = 0.825
$ws.Cells.Item(124,28).Value = -1
$ws.Cells.Item(124,29).Value = 0.875

# Row 178
$ws.Cells.Item(178,2).Value = 7823446
$ws.Cells.Item(178,6).Value = "CS Constantine"
$ws.Cells.Item(178,7).Value = "ASO Chlef"
$ws.Cells.Item(178,8).Value = 3
$ws.Cells.Item(178,9).Value = 1
$ws.Cells.Item(178,10).Value = "H"
$ws.Cells.Item(178,11).Value = 1.444
$ws.Cells.Item(178,12).Value = 3.5
$ws.Cells.Item(178,13).Value = 7.5
$ws.Cells.Item(178,14).Value = 1.45
$ws.Cells.Item(178,15).Value = 3.6
$ws.Cells.Item(178,16).Value = 7
$ws.Cells.Item(178,17).Value = -1
$ws.Cells.Item(178,18).Value = 1.8
$ws.Cells.Item(178,19).Value = 2
$ws.Cells.Item(178,20).Value = 2.25
$ws.Cells.Item(178,21).Value = 1.925
$ws.Cells.Item(178,22).Value = 1.875
$ws.Cells.Item(178,23).Value = 0.45
$ws.Cells.Item(178,24).Value = -1
$ws.Cells.Item(178,25).Value = -1
$ws.Cells.Item(178,26).Value = 0.8
$ws.Cells.Item(178,27).Value = -1
$ws.Cells.Item(178,28).Value = 0.925
$ws.Cells.Item(178,29).Value = -1

# Row 179
$ws.Cells.Item(179,2).Value = 7823445
$ws.Cells.Item(179,6).Value = "MC Oran"
$ws.Cells.Item(179,7).Value = "US Souf"
$ws.Cells.Item(179,8).Value = 2
$ws.Cells.Item(179,9).Value = 1
$ws.Cells.Item(179,10).Value = "H"
$ws.Cells.Item(179,11).Value = 1.666
$ws.Cells.Item(179,12).Value = 3.25
$ws.Cells.Item(179,13).Value = 5
$ws.Cells.Item(179,14).Value = 1.25
$ws.Cells.Item(179,15).Value = 4.333
$ws.Cells.Item(179,16).Value = 11
$ws.Cells.Item(179,17).Value = -1.5
$ws.Cells.Item(179,18).Value = 1.95
$ws.Cells.Item(179,19).Value = 1.85
$ws.Cells.Item(179,20).Value = 2.25
$ws.Cells.Item(179,21).Value = 1.875
$ws.Cells.Item(179,22).Value = 1.925
$ws.Cells.Item(179,23).Value = 0.25
$ws.Cells.Item(179,24).Value = -1
$ws.Cells.Item(179,25).Value = -1
$ws.Cells.Item(179,26).Value = -1
$ws.Cells.Item(179,27).Value = 0.8500000000000001
$ws.Cells.Item(179,28).Value = 0.875
$ws.Cells.Item(179,29).Value = -1

# Row 199
$ws.Cells.Item(199,2).Value = 7971570
$ws.Cells.Item(199,6).Value = "JS Kabylie"
$ws.Cells.Item(199,7).Value = "CR Belouizdad"
$ws.Cells.Item(199,8).Value = 0
$ws.Cells.Item(199,9).Value = 1
$ws.Cells.Item(199,10).Value = "A"
$ws.Cells.Item(199,11).Value = 2.25
$ws.Cells.Item(199,12).Value = 2.875
$ws.Cells.Item(199,13).Value = 3.1
$ws.Cells.Item(199,14).Value = 2.4
$ws.Cells.Item(199,15).Value = 2.8
$ws.Cells.Item(199,16).Value = 2.8
$ws.Cells.Item(199,17).Value = 0
$ws.Cells.Item(199,18).Value = 1.975
$ws.Cells.Item(199,19).Value = 1.825
$ws.Cells.Item(199,20).Value = 1.75
$ws.Cells.Item(199,21).Value = 1.775
$ws.Cells.Item(199,22).Value = 2.025
$ws.Cells.Item(199,23).Value = -1
$ws.Cells.Item(199,24).Value = -1
$ws.Cells.Item(199,25).Value = 1.8
$ws.Cells.Item(199,26).Value = -1
$ws.Cells.Item(199,27).Value = 0.825
$ws.Cells.Item(199,28).Value = -1
$ws.Cells.Item(199,29).Value = 1.025

# Row 200
$ws.Cells.Item(200,2).Value = 7971568
$ws.Cells.Item(200,6).Value = "ES Setif"
$ws.Cells.Item(200,7).Value = "ASO Chlef"
$ws.Cells.Item(200,8).Value = 0
$ws.Cells.Item(200,9).Value = 0
$ws.Cells.Item(200,10).Value = "D"
$ws.Cells.Item(200,11).Value = 1.5
$ws.Cells.Item(200,12).Value = 3.6
$ws.Cells.Item(200,13).Value = 6
$ws.Cells.Item(200,14).Value = 1.5
$ws.Cells.Item(200,15).Value = 3.8
$ws.Cells.Item(200,16).Value = 5.25
$ws.Cells.Item(200,17).Value = -1
$ws.Cells.Item(200,18).Value = 1.925
$ws.Cells.Item(200,19).Value = 1.875
$ws.Cells.Item(200,20).Value = 2.5
$ws.Cells.Item(200,21).Value = 1.95
$ws.Cells.Item(200,22).Value = 1.85
$ws.Cells.Item(200,23).Value = -1
$ws.Cells.Item(200,24).Value = 2.8
$ws.Cells.Item(200,25).Value = -1
$ws.Cells.Item(200,26).Value = -1
$ws.Cells.Item(200,27).Value = 0.875
$ws.Cells.Item(200,28).Value = -1
$ws.Cells.Item(200,29).Value = 0.8500000000000001

# ---- Row 204 odds update ----
$ws.Cells.Item(204,14).Value = 1.5
$ws.Cells.Item(204,15).Value = 3.5
$ws.Cells.Item(204,16).Value = 7
$ws.Cells.Item(204,17).Value = -1
$ws.Cells.Item(204,18).Value = 1.925
$ws.Cells.Item(204,19).Value = 1.875
$ws.Cells.Item(204,20).Value = 2
$ws.Cells.Item(204,21).Value = 1.875
$ws.Cells.Item(204,22).Value = 1.925
$ws.Cells.Item(204,23).Value = 0
$ws.Cells.Item(204,24).Value = 0
$ws.Cells.Item(204,25).Value = 0
$ws.Cells.Item(204,26).Value = 0
$ws.Cells.Item(204,27).Value = 0

# ---- New rows 205-207 ----
# Row 205: copy id/date cell formatting from row 204
$ws.Cells.Item(204,1).Copy() | Out-Null
$ws.Cells.Item(205,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(204,5).Copy() | Out-Null
$ws.Cells.Item(205,5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(205,1).Value = 203
$ws.Cells.Item(205,2).Value = 8007202
$ws.Cells.Item(205,3).Value = "Algeria Division 1"
$ws.Cells.Item(205,4).Value = "Algeria Division 1"
$ws.Cells.Item(205,5).Value = 45387.48958333334
$ws.Cells.Item(205,6).Value = "MC El Bayadh"
$ws.Cells.Item(205,7).Value = "NC Magra"
$ws.Cells.Item(205,11).Value = 1.533
$ws.Cells.Item(205,12).Value = 3.6
$ws.Cells.Item(205,13).Value = 5.5
$ws.Cells.Item(205,14).Value = 1.571
$ws.Cells.Item(205,15).Value = 3.5
$ws.Cells.Item(205,16).Value = 5.25
$ws.Cells.Item(205,17).Value = -0.75
$ws.Cells.Item(205,18).Value = 1.775
$ws.Cells.Item(205,19).Value = 2.025
$ws.Cells.Item(205,20).Value = 2
$ws.Cells.Item(205,21).Value = 1.825
$ws.Cells.Item(205,22).Value = 1.975
$ws.Cells.Item(205,23).Value = 0
$ws.Cells.Item(205,24).Value = 0
$ws.Cells.Item(205,25).Value = 0
$ws.Cells.Item(205,26).Value = 0
$ws.Cells.Item(205,27).Value = 0

# Row 206: copy id/date cell formatting from row 204
$ws.Cells.Item(204,1).Copy() | Out-Null
$ws.Cells.Item(206,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(204,5).Copy() | Out-Null
$ws.Cells.Item(206,5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(206,1).Value = 204
$ws.Cells.Item(206,2).Value = 8007204
$ws.Cells.Item(206,3).Value = "Algeria Division 1"
$ws.Cells.Item(206,4).Value = "Algeria Division 1"
$ws.Cells.Item(206,5).Value = 45387.75
$ws.Cells.Item(206,6).Value = "US Biskra"
$ws.Cells.Item(206,7).Value = "JS Kabylie"
$ws.Cells.Item(206,11).Value = 2
$ws.Cells.Item(206,12).Value = 3
$ws.Cells.Item(206,13).Value = 3.6
$ws.Cells.Item(206,14).Value = 2.1
$ws.Cells.Item(206,15).Value = 2.9
$ws.Cells.Item(206,16).Value = 3.3
$ws.Cells.Item(206,17).Value = -0.25
$ws.Cells.Item(206,18).Value = 1.875
$ws.Cells.Item(206,19).Value = 1.925
$ws.Cells.Item(206,20).Value = 2
$ws.Cells.Item(206,21).Value = 1.825
$ws.Cells.Item(206,22).Value = 1.975
$ws.Cells.Item(206,23).Value = 0
$ws.Cells.Item(206,24).Value = 0
$ws.Cells.Item(206,25).Value = 0
$ws.Cells.Item(206,26).Value = 0
$ws.Cells.Item(206,27).Value = 0

# Row 207: copy id/date cell formatting from row 204
$ws.Cells.Item(204,1).Copy() | Out-Null
$ws.Cells.Item(207,1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(204,5).Copy() | Out-Null
$ws.Cells.Item(207,5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(207,1).Value = 205
$ws.Cells.Item(207,2).Value = 8007205
$ws.Cells.Item(207,3).Value = "Algeria Division 1"
$ws.Cells.Item(207,4).Value = "Algeria Division 1"
$ws.Cells.Item(207,5).Value = 45387.75
$ws.Cells.Item(207,6).Value = "ASO Chlef"
$ws.Cells.Item(207,7).Value = "JS Saoura"
$ws.Cells.Item(207,11).Value = 1.571
$ws.Cells.Item(207,12).Value = 3.6
$ws.Cells.Item(207,13).Value = 5
$ws.Cells.Item(207,14).Value = 1.85
$ws.Cells.Item(207,15).Value = 3.3
$ws.Cells.Item(207,16).Value = 3.6
$ws.Cells.Item(207,17).Value = -0.5
$ws.Cells.Item(207,18).Value = 1.975
$ws.Cells.Item(207,19).Value = 1.825
$ws.Cells.Item(207,20).Value = 2
$ws.Cells.Item(207,21).Value = 1.825
$ws.Cells.Item(207,22).Value = 1.975
$ws.Cells.Item(207,23).Value = 0
$ws.Cells.Item(207,24).Value = 0
$ws.Cells.Item(207,25).Value = 0
$ws.Cells.Item(207,26).Value = 0
$ws.Cells.Item(207,27).Value = 0

$excel.CutCopyMode = 0

